$d = $word.ActiveDocument

# The document accidentally contains the user story
# "Als klant wil ik een kruimelpadfunctie" twice. Remove the duplicate
# paragraph that lives under the "Vervallen:" heading (styled
# "List Paragraph"), leaving the original occurrence (styled
# "Body Text", earlier in the document) untouched.

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt -like "*Als klant wil ik een kruimelpadfunctie*") {
        $styleName = $p.Range.ParagraphStyle.NameLocal
        if ($styleName -eq "List Paragraph") {
            $p.Range.Delete()
        }
    }
}
